# Renaming settings to use underscores, and simplifying the choice_filter
# formula for the "country" question.

$wb = $excel.ActiveWorkbook

# --- survey sheet: update choice_filter formula for "country" (row 18, col I) ---
$survey = $wb.Worksheets.Item("survey")
$survey.Range("I18").Value = "context.region === data('region')"

# --- settings sheet: rename formId/formVersion/formTitle to snake_case ---
$settings = $wb.Worksheets.Item("settings")
$settings.Range("A2").Value = "form_id"
$settings.Range("A3").Value = "form_version"
$settings.Range("A4").Value = "form_title"
